$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" field text (9/21/2021 -> 9/22/2021)
#    on the slide master and every slide layout's Date placeholder.
# ---------------------------------------------------------------------------
function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "9/21/2021") {
                $tr.Text = "9/22/2021"
            }
        }
    }
}

Update-DateField $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateField $layouts.Item($j).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 1: update the title text.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Demo: Unified Data Catalog for Aircraft Images"

# ---------------------------------------------------------------------------
# 3) Slide 1: update the "Scenario 4: ..." content placeholder, collapsing
#    the previous three runs ("Scenario 4: " / "ContentSearch" / " of Photo
#    Metadata to Enrich Catalog") into a single run of new text.
# ---------------------------------------------------------------------------
$contentShape = $slide1.Shapes.Item(3)
$contentRange = $contentShape.TextFrame.TextRange
$contentRange.Delete()
$contentRange.Text = "Scenario 4: Extract Geospatial Metadata to Enrich Catalog"
